$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Khryz Ervon L. Carreon"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2022491"
$ws.Range("C2").Value = "Present"

$ws.Range("A1").Select()
